$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.568.22'
$ws.Range('E2').Value = '  -3.29%  '
$ws.Range('D3').Value = '3.339.91'
$ws.Range('E3').Value = '  -3.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.88'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.31'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.88'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.413'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D12').Value = '3.908.31'
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.63'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').Value = '3.350.52'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('E16').Value = '  -2.08%  '
$ws.Range('D17').Value = '60.565.23'
$ws.Range('E17').Value = '  -3.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.24'
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.47'
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.88'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.81'
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.558'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.73'
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = '3.487.74'
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('E26').Value = '  -5.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.173'
$ws.Range('E27').Value = '  -4.96%  '
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  -3.91%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.08'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.66'
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.86'
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.30'
$ws.Range('E34').Value = '  -3.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.27'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('E36').Value = '  -5.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.80'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '166.70'
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.91'
$ws.Range('E39').Value = '  -12.56%  '
$ws.Range('D40').Value = '3.376.63'
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0746'
$ws.Range('E41').Value = '  -3.50%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.18'
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.757'
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.28'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.14'
$ws.Range('E45').Value = '  -3.27%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.60'
$ws.Range('E46').Value = '  -5.04%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.449.33'
$ws.Range('E47').Value = '  -4.88%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.67'
$ws.Range('E48').Value = '  -3.29%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.29'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0260'
$ws.Range('E51').Value = '  -2.19%  '
